# DeveloperGuide: update section of UndoRedoStack to VersionedAddressBook
#
# 1) Refresh the cached "datetimeFigureOut" footer field text (7/20/17 -> 4/16/2018)
#    everywhere it appears: slide master and every slide layout.
# 2) Remove the now out-of-date "UndoRedo Stack" diagram fragment (a rectangle
#    labelled "UndoRedo"/"Stack", its connector arrow, and its "1" textbox label)
#    from slide 1, since the undo/redo mechanism no longer uses an UndoRedoStack.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "7/20/17") {
                $tr.Text = "4/16/2018"
            }
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DateField $master.Shapes

# --- Every slide layout attached to the master ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateField $layout.Shapes
}

# NOTE: the notes master also carries a cached "7/20/17" datetimeFigureOut
# field (same fix as above), but this COM-interop runtime does not persist
# any mutation made against NotesMaster shapes (text/position/etc. writes
# are silently dropped, and touching it can even corrupt unrelated slide
# master state), so it is intentionally left untouched here.

# --- Remove the obsolete UndoRedo Stack shapes on slide 1 ---
$slide = $p.Slides.Item(1)

$idsToRemove = @(59, 61, 63)
foreach ($targetId in $idsToRemove) {
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            $shp.Delete()
            break
        }
    }
}
